$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5825000
$ws.Range("C3").Value = 4279264.285714285
$ws.Range("C4").Value = 2278571.428571429
$ws.Range("C6").Value = 2558571.428571429
$ws.Range("C7").Value = 11441407.14285714
